$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '49.912.59'
$ws.Range('E2').Value = '  +3.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.647.62'
$ws.Range('E3').Value = '  +6.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '113.96'
$ws.Range('E5').Value = '  +7.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '326.66'
$ws.Range('E6').Value = '  +2.17%  '
$ws.Range('E7').Value = '  +1.59%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.554'
$ws.Range('E9').Value = '  +3.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.99'
$ws.Range('E10').Value = '  +5.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.13'
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0820'
$ws.Range('E12').Value = '  +2.22%  '
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.35'
$ws.Range('E14').Value = '  +3.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.062.52'
$ws.Range('E15').Value = '  +5.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.643.25'
$ws.Range('E16').Value = '  +5.66%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.872'
$ws.Range('E17').Value = '  +5.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '49.825.65'
$ws.Range('E18').Value = '  +3.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.14'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('E20').Value = '  +2.27%  '
$ws.Range('E21').Value = '  -1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0957'
$ws.Range('E22').Value = '  +2.82%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.05'
$ws.Range('E23').Value = '  +1.37%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '277.17'
$ws.Range('E24').Value = '  +2.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.58'
$ws.Range('E25').Value = '  +2.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.78'
$ws.Range('E26').Value = '  +3.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.98'
$ws.Range('E28').Value = '  +2.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('E29').Value = '  -4.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.08'
$ws.Range('E30').Value = '  +3.65%  '
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.18'
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.45'
$ws.Range('E33').Value = '  +3.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.54'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('E35').Value = '  +4.98%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E37').Value = '  +7.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.86'
$ws.Range('E38').Value = '  +5.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.10'
$ws.Range('E39').Value = '  +8.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '124.18'
$ws.Range('E40').Value = '  +1.78%  '
$ws.Range('E41').Value = '  +1.92%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.04'
$ws.Range('E43').Value = '  -1.93%  '
$ws.Range('E44').Value = '  +4.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.079.93'
$ws.Range('E45').Value = '  +4.12%  '
$ws.Range('E46').Value = '  +5.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.32'
$ws.Range('E47').Value = '  +16.17%  '
$ws.Range('E48').Value = '  +4.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.12'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.38'
$ws.Range('E50').Value = '  +4.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.23'
$ws.Range('E51').Value = '  +4.59%  '
